$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "72.949.31"
$ws.Cells.Item(2, 5).Value = "  +3.05%  "

$ws.Cells.Item(3, 4).Value = "3.978.17"
$ws.Cells.Item(3, 5).Value = "  +1.04%  "

$ws.Cells.Item(4, 5).Value = "  -0.08%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "595.39"
$ws.Cells.Item(5, 5).Value = "  +11.47%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "163.21"
$ws.Cells.Item(6, 5).Value = "  +10.84%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.681"
$ws.Cells.Item(7, 5).Value = "  -0.58%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.999"
$ws.Cells.Item(8, 5).Value = "  -0.09%  "

$ws.Cells.Item(9, 5).Value = "  +1.63%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "54.58"
$ws.Cells.Item(11, 5).Value = "  -0.87%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.0000319"
$ws.Cells.Item(12, 5).Value = "  +1.57%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "10.93"
$ws.Cells.Item(13, 5).Value = "  +3.71%  "

$ws.Cells.Item(14, 4).Value = "4.624.28"
$ws.Cells.Item(14, 5).Value = "  +1.40%  "

$ws.Cells.Item(15, 4).Value = "3.981.78"
$ws.Cells.Item(15, 5).Value = "  +1.12%  "

$ws.Cells.Item(16, 5).Value = "  +9.37%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "14.04"
$ws.Cells.Item(17, 5).Value = "  +1.78%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "20.33"
$ws.Cells.Item(18, 5).Value = "  -0.59%  "

$ws.Cells.Item(19, 5).Value = "  +0.33%  "

$ws.Cells.Item(20, 4).Value = "72.610.31"
$ws.Cells.Item(20, 5).Value = "  +2.66%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "438.62"
$ws.Cells.Item(21, 5).Value = "  +4.25%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "4.73"
$ws.Cells.Item(22, 5).Value = "  +12.70%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "95.96"
$ws.Cells.Item(23, 5).Value = "  -1.01%  "

$ws.Cells.Item(24, 5).Value = "  -4.45%  "

$ws.Cells.Item(25, 5).Value = "  -1.03%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "4.33"
$ws.Cells.Item(26, 5).Value = "  +14.27%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "11.26"
$ws.Cells.Item(27, 5).Value = "  -0.20%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "5.95"
$ws.Cells.Item(28, 5).Value = "  +1.45%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "10.31"
$ws.Cells.Item(29, 5).Value = "  -2.79%  "

$ws.Cells.Item(30, 5).Value = "  -0.01%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "7.78"
$ws.Cells.Item(31, 5).Value = "  -0.27%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "13.73"
$ws.Cells.Item(32, 5).Value = "  +3.35%  "

$ws.Cells.Item(33, 5).Value = "  -0.26%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "48.12"
$ws.Cells.Item(34, 5).Value = "  -5.47%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "668.30"
$ws.Cells.Item(35, 5).Value = "  -2.28%  "

$ws.Cells.Item(36, 5).Value = "  +8.66%  "

$ws.Cells.Item(37, 4).Value = "0.0₃0902"
$ws.Cells.Item(37, 5).Value = "  +11.44%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.435"
$ws.Cells.Item(38, 5).Value = "  -0.60%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.999"
$ws.Cells.Item(39, 5).Value = "  +0.03%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "3.34"
$ws.Cells.Item(40, 5).Value = "  -0.48%  "

$ws.Cells.Item(41, 5).Value = "  +5.16%  "

$ws.Cells.Item(42, 5).Value = "  -2.10%  "

$ws.Cells.Item(43, 5).Value = "  +0.22%  "

$ws.Cells.Item(44, 5).Value = "  +2.05%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "10.59"
$ws.Cells.Item(45, 5).Value = "  +6.56%  "

$ws.Cells.Item(46, 5).Value = "  +0.46%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "3.42"
$ws.Cells.Item(47, 5).Value = "  +2.64%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "2.62"
$ws.Cells.Item(48, 5).Value = "  -1.33%  "

$ws.Cells.Item(49, 4).Value = "2.886.36"
$ws.Cells.Item(49, 5).Value = "  +9.81%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "3.05"
$ws.Cells.Item(50, 5).Value = "  +2.31%  "

$ws.Cells.Item(51, 5).Value = "  +4.49%  "
